# CryCompanywiseStockReport_1.xlsx — stock-quantity/value corrections.
# Columns: A=Sl.No  B=Item Code  C=Item Name  D=Rate  E=MRP  F=Qty  G=Value
# "Sub Total:" rows carry the per-company total in column B;
# row 619/620 carry the overall Sub Total / Grand Total.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 3M INDIA LTD ---------------------------------------------------
$ws.Range("F6").Value  = 67
$ws.Range("G6").Value  = 2001.96
$ws.Range("B10").Value = 27263.23        # Sub Total:

# --- BHAWAR SALES CORPORATION (Gillette Presto row) ------------------
$ws.Range("F77").Value = 246
$ws.Range("G77").Value = 11498.04
$ws.Range("B90").Value = 174086.49       # Sub Total:

# --- COLGATE PALMOLIVE INDIA LTD -------------------------------------
$ws.Range("F115").Value = 197
$ws.Range("G115").Value = 19071.57
$ws.Range("B117").Value = 12929.45       # Sub Total:

# --- DESAI (Ginger Garlic Paste) -------------------------------------
$ws.Range("F135").Value = 22
$ws.Range("G135").Value = 682.66
$ws.Range("B138").Value = 2513.8         # Sub Total:

# --- HIM- (lotion row, quantity-only edit) ---------------------------
$ws.Range("F190").Value = 2
$ws.Range("G190").Value = 164.02

# --- HIM-GENTLE BABY SOAP 75G: rows 192/193 swapped content ----------
$ws.Range("B192").Value = 64973
$ws.Range("E192").Value = 35.4
$ws.Range("F192").Value = 2
$ws.Range("G192").Value = 66.6

$ws.Range("B193").Value = 48706
$ws.Range("E193").Value = 39.8
$ws.Range("F193").Value = -144
$ws.Range("G193").Value = -4795.2

$ws.Range("B216").Value = 37572.68       # Sub Total:

# --- HUL-Kissan nango jam 490g: rows 227/228 swapped content ---------
$ws.Range("B227").Value = 55373
$ws.Range("E227").Value = 163.62
$ws.Range("F227").Value = -94
$ws.Range("G227").Value = -13562.32

$ws.Range("B228").Value = 63520
$ws.Range("E228").Value = 153.4
$ws.Range("F228").Value = 65
$ws.Range("G228").Value = 9378.2

# --- JYT - Ujala IDD Top load Detergent Liquid 2lt --------------------
$ws.Range("F303").Value = 27
$ws.Range("G303").Value = 5694.03
$ws.Range("B304").Value = 171097.65      # Sub Total:

# --- KUS-Floor Wiper: rows 322/323 swapped content --------------------
$ws.Range("B322").Value = 47097
$ws.Range("D322").Value = 112.28
$ws.Range("E322").Value = 134.16
$ws.Range("F322").Value = 15
$ws.Range("G322").Value = 1684.2

$ws.Range("B323").Value = 58047
$ws.Range("D323").Value = 105.54
$ws.Range("E323").Value = 126.1
$ws.Range("F323").Value = 39
$ws.Range("G323").Value = 4116.06

# --- CRE-Butter cremfills 100gm: rows 366/367 swapped content ---------
$ws.Range("B366").Value = 53263
$ws.Range("E366").Value = 15.29
$ws.Range("F366").Value = -309
$ws.Range("G366").Value = -3958.29

$ws.Range("B367").Value = 65066
$ws.Range("E367").Value = 13.61
$ws.Range("F367").Value = 90
$ws.Range("G367").Value = 1152.9

# --- CRE-Cremica Honey Oatmeal Cookies: rows 375/376 swapped ----------
$ws.Range("B375").Value = 64927
$ws.Range("E375").Value = 17.26
$ws.Range("F375").Value = 106
$ws.Range("G375").Value = 1719.32

$ws.Range("B376").Value = 45718
$ws.Range("E376").Value = 19.38
$ws.Range("F376").Value = -294
$ws.Range("G376").Value = -4768.68

# --- CRE-Cremica Oatmeal Digestive 112.5 Gm: rows 380/381 swapped -----
$ws.Range("B380").Value = 45709
$ws.Range("E380").Value = 15.69
$ws.Range("F380").Value = -300
$ws.Range("G380").Value = -3945

$ws.Range("B381").Value = 64925
$ws.Range("E381").Value = 13.97
$ws.Range("F381").Value = 111
$ws.Range("G381").Value = 1459.65

# --- PRI-B-50 VIMAL Copper Glass: rows 442/443 swapped ----------------
$ws.Range("B442").Value = 53319
$ws.Range("E442").Value = 310.64
$ws.Range("F442").Value = -6
$ws.Range("G442").Value = -1643.52

$ws.Range("B443").Value = 64810
$ws.Range("E443").Value = 291.22
$ws.Range("F443").Value = 4
$ws.Range("G443").Value = 1095.68

# --- Rasna Nagpur Orange (32 Glass): rows 473/474 swapped -------------
$ws.Range("B473").Value = 60022
$ws.Range("E473").Value = 37.22
$ws.Range("F473").Value = -113
$ws.Range("G473").Value = -3709.79

$ws.Range("B474").Value = 64830
$ws.Range("E474").Value = 34.9
$ws.Range("F474").Value = 107
$ws.Range("G474").Value = 3512.81

# --- Shankys Tip Top Hing Jeera Peanut: rows 572/573 swapped ----------
$ws.Range("B572").Value = 65362
$ws.Range("F572").Value = 20
$ws.Range("G572").Value = 817.4

$ws.Range("B573").Value = 65079
$ws.Range("F573").Value = 6
$ws.Range("G573").Value = 245.22

# --- VVD Priyam Cold Pressed Groundnut Oil Pouch 1 Ltr ----------------
$ws.Range("F599").Value = 1499
$ws.Range("G599").Value = 244501.89
$ws.Range("B606").Value = 401549.01      # Sub Total:

# --- Overall Sub Total / Grand Total ----------------------------------
$ws.Range("B619").Value = 1673811.32     # Sub Total:
$ws.Range("B620").Value = 1673811.32     # Grand Total:
